$wb = $excel.ActiveWorkbook

$wsA = $wb.Worksheets.Item("ALC")
$wsB = $wb.Worksheets.Item("ARM")
$wsC = $wb.Worksheets.Item("BSM")
$wsD = $wb.Worksheets.Item("CUL")
$wsE = $wb.Worksheets.Item("GSM")
$wsF = $wb.Worksheets.Item("LTW")
$wsG = $wb.Worksheets.Item("WVR")

# ALC row 100
$wsA.Range("H100").Value = 3073.8572
$wsA.Range("J100").Value = 3203.4
$wsA.Range("L100").Value = 3203.4
$wsA.Range("N100").Value = -4285.4

# ALC row 129
$wsA.Range("H129").Value = 1280.2234
$wsA.Range("I129").Value = 419.25
$wsA.Range("J129").Value = 1406.2195
$wsA.Range("K129").Value = 1257.75
$wsA.Range("L129").Value = 4218.6585
$wsA.Range("M129").Value = 3742.25
$wsA.Range("N129").Value = -14218.6585

# ALC row 132
$wsA.Range("H132").Value = 4273.686
$wsA.Range("I132").Value = 2306.1904
$wsA.Range("J132").Value = 13455.333
$wsA.Range("K132").Value = 6918.5712
$wsA.Range("L132").Value = 40365.999
$wsA.Range("M132").Value = -4388.5712
$wsA.Range("N132").Value = -45425.999

# ARM row 32
$wsB.Range("H32").Value = 5256
$wsB.Range("I32").Value = 5325.9717
$wsB.Range("J32").Value = 5007.6
$wsB.Range("K32").Value = 5325.9717
$wsB.Range("L32").Value = 5007.6
$wsB.Range("M32").Value = -5038.9717
$wsB.Range("N32").Value = -5581.6

# ARM row 74
$wsB.Range("H74").Value = 1856.9714
$wsB.Range("I74").Value = 1037.6296
$wsB.Range("J74").Value = 4622.25
$wsB.Range("K74").Value = 1037.6296
$wsB.Range("L74").Value = 4622.25
$wsB.Range("M74").Value = -163.6296
$wsB.Range("N74").Value = -6370.25

# ARM row 77
$wsB.Range("H77").Value = 1856.9714
$wsB.Range("I77").Value = 1037.6296
$wsB.Range("J77").Value = 4622.25
$wsB.Range("K77").Value = 5188.148
$wsB.Range("L77").Value = 23111.25
$wsB.Range("M77").Value = -820.1480000000001
$wsB.Range("N77").Value = -31847.25

# ARM row 102
$wsB.Range("H102").Value = 4120.76
$wsB.Range("I102").Value = 4565.95
$wsB.Range("J102").Value = 2340
$wsB.Range("K102").Value = 4565.95
$wsB.Range("L102").Value = 2340
$wsB.Range("M102").Value = -2943.95
$wsB.Range("N102").Value = -5584

# BSM row 22
$wsC.Range("H22").Value = 77.14286
$wsC.Range("I22").Value = 85
$wsC.Range("J22").Value = 66.666664
$wsC.Range("K22").Value = 85
$wsC.Range("L22").Value = 66.666664
$wsC.Range("M22").Value = 88
$wsC.Range("N22").Value = -412.666664

# BSM row 97
$wsC.Range("H97").Value = 15052.857
$wsC.Range("I97").Value = 6979.8
$wsC.Range("J97").Value = 35235.5
$wsC.Range("K97").Value = 6979.8
$wsC.Range("L97").Value = 35235.5
$wsC.Range("M97").Value = -5988.8
$wsC.Range("N97").Value = -37217.5

# BSM row 134
$wsC.Range("H134").Value = 3852.457
$wsC.Range("I134").Value = 2556.3704
$wsC.Range("J134").Value = 8226.75
$wsC.Range("K134").Value = 7669.111199999999
$wsC.Range("L134").Value = 24680.25
$wsC.Range("M134").Value = -5134.111199999999
$wsC.Range("N134").Value = -29750.25

# CUL row 100
$wsD.Range("H100").Value = 4478.75
$wsD.Range("J100").Value = 4478.75
$wsD.Range("L100").Value = 13436.25
$wsD.Range("N100").Value = -15058.25

# CUL row 131
$wsD.Range("H131").Value = 860.61
$wsD.Range("I131").Value = 389.73685
$wsD.Range("J131").Value = 971.0617
$wsD.Range("K131").Value = 1169.21055
$wsD.Range("L131").Value = 2913.1851
$wsD.Range("M131").Value = 3870.78945
$wsD.Range("N131").Value = -12993.1851

# CUL row 132
$wsD.Range("H132").Value = 2774.9546
$wsD.Range("I132").Value = 804.9091
$wsD.Range("J132").Value = 4745
$wsD.Range("K132").Value = 7244.1819
$wsD.Range("L132").Value = 42705
$wsD.Range("M132").Value = -4714.1819
$wsD.Range("N132").Value = -47765

# GSM row 124
$wsE.Range("H124").Value = 38121.21
$wsE.Range("I124").Value = 25000
$wsE.Range("J124").Value = 38531.25
$wsE.Range("K124").Value = 25000
$wsE.Range("L124").Value = 38531.25
$wsE.Range("M124").Value = -20090
$wsE.Range("N124").Value = -48351.25

# LTW row 16
$wsF.Range("H16").Value = 2634.7727
$wsF.Range("I16").Value = 1956.3889
$wsF.Range("K16").Value = 1956.3889
$wsF.Range("M16").Value = -1786.3889

# LTW row 22
$wsF.Range("H22").Value = 650.5
$wsF.Range("I22").Value = 700
$wsF.Range("J22").Value = 601
$wsF.Range("K22").Value = 700
$wsF.Range("L22").Value = 601
$wsF.Range("M22").Value = -405
$wsF.Range("N22").Value = -1191

# LTW row 27
$wsF.Range("H27").Value = 650.5
$wsF.Range("I27").Value = 700
$wsF.Range("J27").Value = 601
$wsF.Range("K27").Value = 700
$wsF.Range("L27").Value = 601
$wsF.Range("M27").Value = -593
$wsF.Range("N27").Value = -815

# LTW row 46
$wsF.Range("H46").Value = 1355.5555
$wsF.Range("I46").Value = 1066.6666
$wsF.Range("J46").Value = 1933.3334
$wsF.Range("K46").Value = 1066.6666
$wsF.Range("L46").Value = 1933.3334
$wsF.Range("M46").Value = -878.6666
$wsF.Range("N46").Value = -2309.3334

# LTW row 55
$wsF.Range("H55").Value = 539.41174
$wsF.Range("I55").Value = 274.44446
$wsF.Range("J55").Value = 837.5
$wsF.Range("K55").Value = 274.44446
$wsF.Range("L55").Value = 837.5
$wsF.Range("M55").Value = -101.44446
$wsF.Range("N55").Value = -1183.5

# LTW row 81
$wsF.Range("H81").Value = 38890.5
$wsF.Range("I81").Value = 0
$wsF.Range("J81").Value = 38890.5
$wsF.Range("K81").Value = 0
$wsF.Range("L81").Value = 38890.5
$wsF.Range("N81").Value = -40886.5
$wsF.Range("M81").ClearContents()

# LTW row 82
$wsF.Range("H82").Value = 2065
$wsF.Range("J82").Value = 2137.5
$wsF.Range("L82").Value = 2137.5
$wsF.Range("N82").Value = -2859.5

# LTW row 84
$wsF.Range("H84").Value = 38890.5
$wsF.Range("I84").Value = 0
$wsF.Range("J84").Value = 38890.5
$wsF.Range("K84").Value = 0
$wsF.Range("L84").Value = 116671.5
$wsF.Range("N84").Value = -126655.5
$wsF.Range("M84").ClearContents()

# LTW row 85
$wsF.Range("H85").Value = 2065
$wsF.Range("J85").Value = 2137.5
$wsF.Range("L85").Value = 2137.5
$wsF.Range("N85").Value = -4633.5

# WVR row 100
$wsG.Range("H100").Value = 833.6111
$wsG.Range("I100").Value = 500.5
$wsG.Range("J100").Value = 928.7857
$wsG.Range("K100").Value = 1001
$wsG.Range("L100").Value = 1857.5714
$wsG.Range("M100").Value = -460
$wsG.Range("N100").Value = -2939.5714

# WVR row 132
$wsG.Range("H132").Value = 4739.129
$wsG.Range("I132").Value = 4824.6206
$wsG.Range("J132").Value = 3499.5
$wsG.Range("K132").Value = 14473.8618
$wsG.Range("L132").Value = 10498.5
$wsG.Range("M132").Value = -11943.8618
$wsG.Range("N132").Value = -15558.5
